$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the "Änderungsnummer" (change number) for rows 16-20 ---
$ws.Range("A16").Value = 1015
$ws.Range("A17").Value = 1016
$ws.Range("A18").Value = 1017
$ws.Range("A19").Value = 1018
$ws.Range("A20").Value = 1019

# --- Row 20: "gefunden in Version" moves from 2.28 to 2.29, comment cleared ---
$ws.Range("F20").Value = "2.29"
$ws.Range("G20").ClearContents()

# --- New rows 21-24: copy formatting from row 20 as a starting template ---
$ws.Range("A20:H20").Copy()
$ws.Range("A21:H24").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Row 21
$ws.Range("A21").Value = 1020
$ws.Range("B21").Value = 43952
$ws.Range("C21").Value = "Performance very low"
$ws.Range("D21").Value = "Harold"
$ws.Range("E21").Value = "Fehler"
$ws.Range("F21").Value = "2.29"
$ws.Range("G21").Value = "delay between ARDUINO message too high"
$ws.Range("H21").Value = "2.30"

# Row 22
$ws.Range("A22").Value = 1021
$ws.Range("B22").Value = 43952
$ws.Range("C22").Value = "Allow entry of individual portname (for MAC and LINUX users)"
$ws.Range("D22").Value = "Harold"
$ws.Range("E22").Value = "Neue Funktion"
$ws.Range("H22").Value = "2.30"

# Row 23
$ws.Range("A23").Value = 1032
$ws.Range("B23").Value = 43952
$ws.Range("C23").Value = "Improve group info handling"
$ws.Range("D23").Value = "Harold"
$ws.Range("E23").Value = "Fehler"
$ws.Range("F23").Value = "2.29"
$ws.Range("H23").Value = "2.30"

# Row 24
$ws.Range("A24").Value = 1033
$ws.Range("B24").Value = 43952
$ws.Range("C24").Value = "Springe zur zugehörigen Macroseite beim Acklicken einer LED in der LEDListe"
$ws.Range("D24").Value = "Harold"
$ws.Range("E24").Value = "Neue Funktion"
$ws.Range("H24").Value = "2.30"

# Rows 21 and 24 wrap onto two lines (like rows 1, 14, 17 etc.), so they
# use the taller 30pt row height used elsewhere in the sheet for 2-line rows
$ws.Rows("21").RowHeight = 30
$ws.Rows("24").RowHeight = 30

# --- View state: selection moved to C26 (sheet scrolled so row 7 is on top) ---
$ws.Range("C26").Select()
try { $excel.ActiveWindow.ScrollRow = 7 } catch {}
